# ODF Template Demo 1 (29 04 25).
# Adds four new Q&A rows (320-323) to the query_responses sheet, matching the
# existing "Question / Model Name / Response" table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A320").Value = "How many types of information can an ODT file hold?"
$ws.Range("B320").Value = "llama3.2:latest"
$ws.Range("C320").Value = "According to the provided documents, an ODT file can hold encapsulated geological information such as Gas data, descriptions, lithology, % lithology, table/curve parameter details and qualitative information."

$ws.Range("A321").Value = "Which type of content can be viewed from an ODT file?"
$ws.Range("B321").Value = "llama3.2:latest"
$ws.Range("C321").Value = "According to the Documents section, an ODT file contains:`n1. Library information (headers, lithology, modifiers, structures, and symbols)`n2. View file contents (track layout information, depth and screen units, scale and pen information)`n3. ini file settings (curve defaults, computed curves and table definitions)`nTherefore, an ODT file can be viewed as containing all of these types of content."

$ws.Range("A322").Value = "List the library informations from an ODT file"
$ws.Range("B322").Value = "llama3.2:latest"
$ws.Range("C322").Value = "Based on the provided documents, to list library information from an ODT file, you can:`n1. Open the template by selecting `"Templat`" on the GEOMenu.`n2. In the Open Template dialog box, browse to the location of your ODT file and select it.`n3. Click `"Open`".`nThe document information tree will be displayed, showing warnings if some components are not functioning correctly or have settings that are deemed not best practice."

$ws.Range("A323").Value = "Can an ODT be created from an existing ODF file?"
$ws.Range("B323").Value = "llama3.2:latest"
$ws.Range("C323").Value = "Yes, an ODT (OpenDocument Template) can be created from an existing ODF (OpenDocument Format) file by selecting `"Save as`" and then choosing the template option. Alternatively, you can also create a new ODT directly from an ODF using the `"Create Template`" function in the GEO menu."
